$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75, shifting existing rows 75-80 down to 76-81
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new weekly record
$ws.Cells.Item(75, 1).Value = 7
$ws.Cells.Item(75, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(75, 3).Value = "Ñuble"
$ws.Cells.Item(75, 4).Value = 44858
$ws.Cells.Item(75, 5).Value = 16
$ws.Cells.Item(75, 6).Value = 100112022
$ws.Cells.Item(75, 7).Value = "Arveja Verde"
$ws.Cells.Item(75, 8).Value = "Sin especificar"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 100
$ws.Cells.Item(75, 11).Value = 15000
$ws.Cells.Item(75, 12).Value = 16000
$ws.Cells.Item(75, 13).Value = 15500
$ws.Cells.Item(75, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(75, 15).Value = "Región del Maule"
$ws.Cells.Item(75, 16).Value = 620
$ws.Cells.Item(75, 17).Value = 25
$ws.Cells.Item(75, 18).Value = "Hortaliza"
